$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3786642551422119
$ws.Range("B1").Value = 1.850308656692505
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.568034887313843
$ws.Range("E1").Value = 1.378679513931274
